$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.685.91"
$ws.Range("E2").Value = "  -0.79%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.789.27"
$ws.Range("E3").Value = "  +0.35%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.96"
$ws.Range("E5").Value = "  +0.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.88"
$ws.Range("E6").Value = "  -0.39%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.785.36"
$ws.Range("E7").Value = "  +0.18%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("E9").Value = "  +0.11%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.160"
$ws.Range("E10").Value = "  -0.14%  "

$ws.Range("E11").Value = "  -1.10%  "

$ws.Range("E12").Value = "  -0.23%  "

$ws.Range("E13").Value = "  -2.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.03"
$ws.Range("E14").Value = "  -0.19%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.425.50"
$ws.Range("E15").Value = "  +0.52%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.800.24"
$ws.Range("E16").Value = "  +1.77%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.55"
$ws.Range("E17").Value = "  +3.88%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.653.03"
$ws.Range("E18").Value = "  -0.74%  "

$ws.Range("E19").Value = "  +0.85%  "

$ws.Range("E20").Value = "  +0.09%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.99"
$ws.Range("E21").Value = "  -7.45%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "459.27"
$ws.Range("E22").Value = "  -1.34%  "

$ws.Range("E23").Value = "  +0.14%  "

$ws.Range("E24").Value = "  +3.57%  "

$ws.Range("E25").Value = "  -0.44%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.15"
$ws.Range("E26").Value = "  +2.69%  "

$ws.Range("E27").Value = "  -3.32%  "

$ws.Range("E28").Value = "  -0.03%  "

$ws.Range("E29").Value = "  -1.30%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.939.96"
$ws.Range("E30").Value = "  +0.51%  "

$ws.Range("E31").Value = "  -0.19%  "

$ws.Range("E32").Value = "  +4.26%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.20"
$ws.Range("E33").Value = "  -1.26%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.64"
$ws.Range("E34").Value = "  -0.90%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.12%  "

$ws.Range("E36").Value = "  -0.71%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0998"
$ws.Range("E37").Value = "  -0.58%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.35"
$ws.Range("E38").Value = "  -2.52%  "

$ws.Range("E39").Value = "  -0.40%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.992"
$ws.Range("E40").Value = "  -0.83%  "

$ws.Range("E41").Value = "  -0.11%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.09%  "

$ws.Range("E43").Value = "  -0.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.02"
$ws.Range("E44").Value = "  +2.25%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "43.85"
$ws.Range("E45").Value = "  -1.72%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.297"
$ws.Range("E46").Value = "  -0.91%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "149.82"
$ws.Range("E47").Value = "  +2.87%  "

$ws.Range("E48").Value = "  -1.64%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "26.95"
$ws.Range("E49").Value = "  +6.94%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "390.03"
$ws.Range("E50").Value = "  +0.01%  "

$ws.Range("E51").Value = "  -4.31%  "
